# Update cryptos list values to the latest scraped snapshot (GitHub Actions data refresh).
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h)
# Numeric-looking Price strings are written via a Text-number-format round trip so
# they remain plain text cells (matching the source data's inline-string format)
# instead of being auto-converted to numbers by Excel, then the style is reset back
# to Normal so no stray cell formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.105.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.46%  "

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.828.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.45%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.45%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.41%  "

# Row 6: USDC
$ws.Range("E6").Value = "  -0.39%  "

# Row 7: XRP
$ws.Range("E7").Value = "  +7.37%  "

# Row 8: Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3742"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.00%  "

# Row 9: Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07319"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.33%  "

# Row 10: Polygon
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8629"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.33%  "

# Row 11: Solana
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.96"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.04%  "

# Row 12: WrappedEther
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.824.58"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.15%  "

# Row 13: Chainlink
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.724"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.92%  "

# Row 14: Litecoin
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.63%  "

# Row 16: TRON
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07090"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.32%  "

# Row 17: BinanceUSD
$ws.Range("E17").Value = "  -0.43%  "

# Row 18: ShibaInu
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008846"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.11%  "

# Row 19: Dai
$ws.Range("E19").Value = "  -0.32%  "

# Row 20: Avalanche
$ws.Range("E20").Value = "  +0.01%  "

# Row 21: WrappedBTC
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.106.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.56%  "

# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.203"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.19%  "

# Row 23: Cosmos
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.32%  "

# Row 24: Toncoin
$ws.Range("E24").Value = "  +0.03%  "

# Row 25: Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.73%  "

# Row 26: LidoDAOToken
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.227"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.31%  "

# Row 27: EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.85%  "

# Row 28: InternetComputer(DFINITY)
$ws.Range("E28").Value = "  +1.38%  "

# Row 29: BitcoinCash
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.32%  "

# Row 30: Stellar
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08912"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.28%  "

# Row 31: ImmutableX -> ARBITRUM (row content swap)
$ws.Range("B31").Value = "ARBITRUM"
$ws.Range("C31").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.200"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.66%  "

# Row 32: ARBITRUM -> ImmutableX (row content swap)
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7645"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.08%  "

# Row 33: HuobiToken
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.973"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.37%  "

# Row 34: Filecoin
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.481"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.95%  "

# Row 35: Frax
$ws.Range("E35").Value = "  -0.37%  "

# Row 36: TrustWalletToken
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.105"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.77%  "

# Row 37: VeChain
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01973"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.58%  "

# Row 38: Hedera
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05298"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.65%  "

# Row 39: TheSandbox
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5380"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.16%  "

# Row 40: FraxShare
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.195"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.00%  "

# Row 41: MXToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.884"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.75%  "

# Row 42: Algorand
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1721"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.71%  "

# Row 43: Decentraland
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5221"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +11.15%  "

# Row 44: Aptos
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.652"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.66%  "

# Row 45: EnergySwap
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.57%  "

# Row 46: RenderToken
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.993"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +11.15%  "

# Row 47: Quant
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "106.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.50%  "

# Row 48: NEARProtocol
$ws.Range("E48").Value = "  +1.73%  "

# Row 49: Cronos
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06433"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.75%  "

# Row 50: PaxDollar
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.000"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.43%  "

# Row 51: ThetaToken
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9250"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.44%  "
